$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update email addresses (test.f -> test.com)
$ws.Range("C2").Value = "alexde@test.com"
$ws.Range("C3").Value = "canbartu@test.com"
$ws.Range("C4").Value = "senturks@test.com"

# Remove the now-unused F column data (F2:F4 held duplicate id values)
$ws.Range("F2:F4").Clear()

# Update the active selection as recorded in the sheet view
$ws.Range("B6").Select()
